$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.274.58"
$ws.Range("E2").Value = "  -0.14%  "

$ws.Range("D3").Value = "1.829.33"
$ws.Range("E3").Value = "  -0.58%  "

$ws.Range("E4").Value = "  +0.35%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.69"
$ws.Range("E5").Value = "  -1.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6020"
$ws.Range("E6").Value = "  -3.00%  "

$ws.Range("E7").Value = "  +0.38%  "

$ws.Range("E8").Value = "  -4.67%  "

$ws.Range("E9").Value = "  -3.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.55"
$ws.Range("E10").Value = "  -4.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07610"
$ws.Range("E11").Value = "  -1.50%  "

$ws.Range("D12").Value = "1.872.14"
$ws.Range("E12").Value = "  +2.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.751"
$ws.Range("E13").Value = "  -3.61%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6321"
$ws.Range("E14").Value = "  -3.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009882"
$ws.Range("E15").Value = "  -1.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.75"
$ws.Range("E16").Value = "  -4.35%  "

$ws.Range("D17").Value = "28.956.90"
$ws.Range("E17").Value = "  -1.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.581"
$ws.Range("E18").Value = "  -10.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.32"
$ws.Range("E19").Value = "  -7.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.005"
$ws.Range("E20").Value = "  +0.39%  "

$ws.Range("E21").Value = "  -4.67%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.891"
$ws.Range("E22").Value = "  -3.88%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.007"
$ws.Range("E23").Value = "  +0.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "156.37"
$ws.Range("E24").Value = "  -0.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.978"
$ws.Range("E25").Value = "  -4.71%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1289"
$ws.Range("E26").Value = "  -2.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.53"
$ws.Range("E27").Value = "  -3.71%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06462"
$ws.Range("E28").Value = "  -6.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.418"
$ws.Range("E29").Value = "  -3.86%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.445"
$ws.Range("E30").Value = "  -2.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.834"
$ws.Range("E31").Value = "  -2.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.797"
$ws.Range("E32").Value = "  -5.21%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.097"
$ws.Range("E33").Value = "  -4.81%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.728"
$ws.Range("E34").Value = "  -0.73%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6486"
$ws.Range("E35").Value = "  -4.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.547"
$ws.Range("E36").Value = "  -1.34%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.757"
$ws.Range("E37").Value = "  -0.67%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01757"
$ws.Range("E38").Value = "  -3.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.615"
$ws.Range("E39").Value = "  -0.38%  "

$ws.Range("D40").Value = "1.141.51"
$ws.Range("E40").Value = "  -7.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8934"
$ws.Range("E41").Value = "  -5.39%  "

$ws.Range("E42").Value = "  +0.36%  "

$ws.Range("D43").Value = "1.998.58"
$ws.Range("E43").Value = "  +0.48%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.83"
$ws.Range("E44").Value = "  -0.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.27"
$ws.Range("E45").Value = "  -4.13%  "

$ws.Range("E46").Value = "  -3.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.626"
$ws.Range("E47").Value = "  -3.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.490"
$ws.Range("E48").Value = "  -2.84%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05498"
$ws.Range("E49").Value = "  -2.46%  "

$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4540"
$ws.Range("E50").Value = "  -0.80%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.399"
$ws.Range("E51").Value = "  -6.57%  "

